$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "xLHlk358"
$ws.Range("B2").Value = 23091444
$ws.Range("C2").Value = "chqahhv29"
$ws.Range("D2").Value = "aG7M#2!t"
$ws.Range("F2").Value = "LuvzfEwp"
$ws.Range("G2").Value = "FudG"
